$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '70.913.96'
$ws.Range("E2").Value = '  +0.13%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.848.22'
$ws.Range("E3").Value = '  +1.39%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '708.99'
$ws.Range("E5").Value = '  +0.96%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.64'
$ws.Range("E6").Value = '  -0.18%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.846.50'
$ws.Range("E7").Value = '  +1.32%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("E9").Value = '  -0.49%  '
$ws.Range("E10").Value = '  -0.24%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.28'
$ws.Range("E11").Value = '  -1.52%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.458'
$ws.Range("E12").Value = '  -0.46%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000257'
$ws.Range("E13").Value = '  -0.45%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.71'
$ws.Range("E14").Value = '  +0.90%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.494.71'
$ws.Range("E15").Value = '  +1.34%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.860.07'
$ws.Range("E16").Value = '  +1.70%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '70.966.71'
$ws.Range("E17").Value = '  +0.22%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.19'
$ws.Range("E18").Value = '  +0.03%  '
$ws.Range("E19").Value = '  +0.94%  '
$ws.Range("E20").Value = '  -2.81%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.67'
$ws.Range("E21").Value = '  -4.01%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '492.79'
$ws.Range("E22").Value = '  +2.27%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.717'
$ws.Range("E23").Value = '  +0.32%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '85.20'
$ws.Range("E24").Value = '  +0.70%  '
$ws.Range("E25").Value = '  +2.45%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.61'
$ws.Range("E26").Value = '  +1.59%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.15'
$ws.Range("E27").Value = '  -1.97%  '
$ws.Range("E28").Value = '  -3.03%  '
$ws.Range("E29").Value = '  +1.42%  '
$ws.Range("E30").Value = '  +0.03%  '
$ws.Range("E31").Value = '  -0.32%  '
$ws.Range("E32").Value = '  -0.53%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '29.44'
$ws.Range("E33").Value = '  -0.28%  '
$ws.Range("E34").Value = '  -1.28%  '
$ws.Range("E35").Value = '  -0.63%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.803.97'
$ws.Range("E36").Value = '  +1.58%  '
$ws.Range("E37").Value = '  -0.24%  '
$ws.Range("E38").Value = '  +0.30%  '
$ws.Range("E39").Value = '  +6.87%  '
$ws.Range("E41").Value = '  +6.79%  '
$ws.Range("E42").Value = '  -3.07%  '
$ws.Range("E44").Value = '  +0.18%  '
$ws.Range("E45").Value = '  -3.92%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '162.85'
$ws.Range("E46").Value = '  +0.12%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '48.75'
$ws.Range("E47").Value = '  -0.60%  '
$ws.Range("E48").Value = '  +1.73%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '416.29'
$ws.Range("E49").Value = '  +2.11%  '
$ws.Range("E50").Value = '  -1.27%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.62'
$ws.Range("E51").Value = '  +0.68%  '
